$wb = $excel.ActiveWorkbook

# Duplicate the "Hippo" sheet (Excel inserts the copy immediately before the
# sheet it is copied relative to), then rename it and it becomes the new
# first tab: "Don1-Phase1" - a fresh training-stats sheet for the Don1 phase 1
# session.
$hippo = $wb.Worksheets.Item("Hippo")
$hippo.Copy($hippo)
$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "Don1-Phase1"

# Record the Don1 phase 1 run's stats (first attempt block) on the new sheet.
$newSheet.Range("A4").Value = 17
$newSheet.Range("B4").Value = 83
$newSheet.Range("C4").Value = 162
$newSheet.Range("J4").Value = 30

# Second attempt block wasn't run yet this session - clear the copied sample
# data but keep the frames-per-hit tally for it.
$newSheet.Range("A5").ClearContents()
$newSheet.Range("B5").ClearContents()
$newSheet.Range("C5").ClearContents()
$newSheet.Range("J5").Value = 36

# Third attempt block also wasn't run yet - clear its copied sample data too.
$newSheet.Range("A6").ClearContents()
$newSheet.Range("B6").ClearContents()
$newSheet.Range("C6").ClearContents()

# The "Hippo" sheet is no longer the active tab - drop its selection back to
# where it was last left (and it's no longer the focused/selected tab).
$hippoRef = $wb.Worksheets.Item("Hippo")
$hippoRef.Activate()
$hippoRef.Range("E7").Select()

# Make the new sheet the active tab with a single-cell selection.
$newSheet.Activate()
$newSheet.Range("C5").Select()
